$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 8): a time-of-day style number in column A and
# the reused shared string "c" in column B, matching the formatting already
# used by the rest of the data rows.
$ws.Range("A8").Value = 0.5
$ws.Range("A8").NumberFormat = "0.000"
$ws.Range("B8").Value = "c"

# Grow the autofilter so it covers the new row and add the corresponding
# filter value (0.500) alongside the two that were already selected.
$ws.AutoFilterMode = $false
$ws.Range("A1:B8").AutoFilter(1, @("0.046", "0.500", "0.516"), 7) | Out-Null

# The hidden _FilterDatabase defined name must track the autofilter range.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Munka1!`$A`$1:`$B`$8"

# Match the updated selection left behind in the saved file.
$ws.Range("C7").Select() | Out-Null
